$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assignments")

# Delete row 3 (the second data row), then delete column A (group_file) and
# the now-shifted last column (mode).
$ws.Rows.Item(3).Delete()
$ws.Columns.Item(1).Delete()
$ws.Columns.Item(7).Delete()

# Update the remaining row 2 values.
$ws.Range("B2").Value = "\\nashp\DATABUHP\Nam SEO\.KhaiThacShort\Number B\232G.mp4"
$ws.Range("C2").Value = "fasgsagas"

# Column widths.
$ws.Columns.Item(1).ColumnWidth = 19
$ws.Columns.Item(2).ColumnWidth = 59
$ws.Columns.Item(3).ColumnWidth = 11
$ws.Columns.Item(4).ColumnWidth = 13
$ws.Columns.Item(5).ColumnWidth = 14
$ws.Columns.Item(6).ColumnWidth = 14

$ws.Range("A1:F2").AutoFilter()
